$d = $word.ActiveDocument

# The document has one section whose header/footer pair is split across
# two physical parts ("first page" header/footer + "default" header/footer).
# Each part carries one inline picture (the BTEC logo in the headers, the
# Pearson logo in the footers). The logos' display names need to be
# renamed:
#   headers: image1.jpg -> image2.jpg
#   footers: image2.png -> image1.png
#
# WdHeaderFooterIndex: 1 = wdHeaderFooterPrimary (the "default" part),
#                       2 = wdHeaderFooterFirstPage (the "first" part).
# Re-fetch the Section/Header/Footer objects fresh for every single
# assignment so no handle is reused after a prior edit invalidates it.

$d.Sections.Item(1).Headers.Item(1).Range.InlineShapes.Item(1).Name = "image2.jpg"
$d.Sections.Item(1).Headers.Item(2).Range.InlineShapes.Item(1).Name = "image2.jpg"

$d.Sections.Item(1).Footers.Item(1).Range.InlineShapes.Item(1).Name = "image1.png"
$d.Sections.Item(1).Footers.Item(2).Range.InlineShapes.Item(1).Name = "image1.png"
